# Auto-generated Excel COM-interop script implementing the target diff.
$wb = $excel.ActiveWorkbook

# --- 1. Rename the existing "metrics" sheet (sheetId=2) to "worksheet" ---
$wsOld = $wb.Worksheets.Item("metrics")
$wsOld.Name = "worksheet"

# --- 2. Insert a brand-new "metrics" sheet right after "worksheet" ---
$wsMetrics = $wb.Worksheets.Add($null, $wsOld)
$wsMetrics.Name = "metrics"

$wsData = $wb.Worksheets.Item("data")
$ws = $wb.Worksheets.Item("worksheet")

# --- 3. Defined names (ACC, BACC, ... TPR) all pointing at worksheet!$<col>$2 ---
$wb.Names.Add("ACC", "=worksheet!`$P`$2")
$wb.Names.Add("BACC", "=worksheet!`$AA`$2")
$wb.Names.Add("CK", "=worksheet!`$AH`$2")
$wb.Names.Add("CKC", "=worksheet!`$AG`$2")
$wb.Names.Add("CRR", "=worksheet!`$R`$2")
$wb.Names.Add("DR", "=worksheet!`$Q`$2")
$wb.Names.Add("F05_", "=worksheet!`$AF`$2")
$wb.Names.Add("F1_", "=worksheet!`$AD`$2")
$wb.Names.Add("F2_", "=worksheet!`$AE`$2")
$wb.Names.Add("FC", "=worksheet!`$G`$2")
$wb.Names.Add("FDR", "=worksheet!`$T`$2")
$wb.Names.Add("FN", "=worksheet!`$D`$2")
$wb.Names.Add("FNR", "=worksheet!`$M`$2")
$wb.Names.Add("FOR", "=worksheet!`$U`$2")
$wb.Names.Add("FP", "=worksheet!`$C`$2")
$wb.Names.Add("FPR", "=worksheet!`$O`$2")
$wb.Names.Add("GM", "=worksheet!`$AB`$2")
$wb.Names.Add("INFORM", "=worksheet!`$Z`$2")
$wb.Names.Add("LRN", "=worksheet!`$Y`$2")
$wb.Names.Add("LRP", "=worksheet!`$X`$2")
$wb.Names.Add("MARK", "=worksheet!`$AC`$2")
$wb.Names.Add("MCC", "=worksheet!`$AI`$2")
$wb.Names.Add("MCR", "=worksheet!`$W`$2")
$wb.Names.Add("N", "=worksheet!`$K`$2")
$wb.Names.Add("NPV", "=worksheet!`$V`$2")
$wb.Names.Add("ON", "=worksheet!`$I`$2")
$wb.Names.Add("OP", "=worksheet!`$H`$2")
$wb.Names.Add("OR", "=worksheet!`$AJ`$2")
$wb.Names.Add("P", "=worksheet!`$J`$2")
$wb.Names.Add("PPV", "=worksheet!`$S`$2")
$wb.Names.Add("PREV", "=worksheet!`$AL`$2")
$wb.Names.Add("SKEW", "=worksheet!`$AM`$2")
$wb.Names.Add("Sn", "=worksheet!`$A`$2")
$wb.Names.Add("TC", "=worksheet!`$F`$2")
$wb.Names.Add("TN", "=worksheet!`$E`$2")
$wb.Names.Add("TNR", "=worksheet!`$N`$2")
$wb.Names.Add("TP", "=worksheet!`$B`$2")
$wb.Names.Add("TPR", "=worksheet!`$L`$2")

# --- 4. Header row (row 1) on "worksheet" ---
$ws.Cells.Item(1, 1).Value = "Sn"
$ws.Cells.Item(1, 2).Value = "TP"
$ws.Cells.Item(1, 3).Value = "FP"
$ws.Cells.Item(1, 4).Value = "FN"
$ws.Cells.Item(1, 5).Value = "TN"
$ws.Cells.Item(1, 6).Value = "TC"
$ws.Cells.Item(1, 7).Value = "FC"
$ws.Cells.Item(1, 8).Value = "OP"
$ws.Cells.Item(1, 9).Value = "ON"
$ws.Cells.Item(1, 10).Value = "P"
$ws.Cells.Item(1, 11).Value = "N"
$ws.Cells.Item(1, 12).Value = "TPR"
$ws.Cells.Item(1, 13).Value = "FNR"
$ws.Cells.Item(1, 14).Value = "TNR"
$ws.Cells.Item(1, 15).Value = "FPR"
$ws.Cells.Item(1, 16).Value = "ACC"
$ws.Cells.Item(1, 17).Value = "DR"
$ws.Cells.Item(1, 18).Value = "CRR"
$ws.Cells.Item(1, 19).Value = "PPV"
$ws.Cells.Item(1, 20).Value = "FDR"
$ws.Cells.Item(1, 21).Value = "FOR"
$ws.Cells.Item(1, 22).Value = "NPV"
$ws.Cells.Item(1, 23).Value = "MCR"
$ws.Cells.Item(1, 24).Value = "LRP"
$ws.Cells.Item(1, 25).Value = "LRN"
$ws.Cells.Item(1, 26).Value = "INFORM"
$ws.Cells.Item(1, 27).Value = "BACC"
$ws.Cells.Item(1, 28).Value = "GM"
$ws.Cells.Item(1, 29).Value = "MARK"
$ws.Cells.Item(1, 30).Value = "F1"
$ws.Cells.Item(1, 31).Value = "F2"
$ws.Cells.Item(1, 32).Value = "F05"
$ws.Cells.Item(1, 33).Value = "CKC"
$ws.Cells.Item(1, 34).Value = "CK"
$ws.Cells.Item(1, 35).Value = "MCC"
$ws.Cells.Item(1, 36).Value = "OR"
$ws.Cells.Item(1, 37).Value = "DP"
$ws.Cells.Item(1, 38).Value = "PREV"
$ws.Cells.Item(1, 39).Value = "SKEW"

# --- 5. Formula row (row 2) on "worksheet" ---
$ws.Cells.Item(2, 1).Formula = "=COUNT(y)"
$ws.Cells.Item(2, 2).Formula = "=COUNTIFS(y,1,y_pred,1)"
$ws.Cells.Item(2, 3).Formula = "=COUNTIFS(y,0,y_pred,1)"
$ws.Cells.Item(2, 4).Formula = "=COUNTIFS(y,1,y_pred,0)"
$ws.Cells.Item(2, 5).Formula = "=COUNTIFS(y,0,y_pred,0)"
$ws.Cells.Item(2, 6).Formula = "=TP+TN"
$ws.Cells.Item(2, 7).Formula = "=FP+FN"
$ws.Cells.Item(2, 8).Formula = "=TP+FP"
$ws.Cells.Item(2, 9).Formula = "=FN+TN"
$ws.Cells.Item(2, 10).Formula = "=TP+FN"
$ws.Cells.Item(2, 11).Formula = "=FP+TN"
$ws.Cells.Item(2, 12).Formula = "=TP/P"
$ws.Cells.Item(2, 13).Formula = "=FN/P"
$ws.Cells.Item(2, 14).Formula = "=TN/N"
$ws.Cells.Item(2, 15).Formula = "=FP/N"
$ws.Cells.Item(2, 16).Formula = "=TC/Sn"
$ws.Cells.Item(2, 17).Formula = "=TP/Sn"
$ws.Cells.Item(2, 18).Formula = "=TN/Sn"
$ws.Cells.Item(2, 19).Formula = "=TP/OP"
$ws.Cells.Item(2, 20).Formula = "=FP/OP"
$ws.Cells.Item(2, 21).Formula = "=FN/ON"
$ws.Cells.Item(2, 22).Formula = "=TN/ON"
$ws.Cells.Item(2, 23).Formula = "=FC/Sn"
$ws.Cells.Item(2, 24).Formula = "=TPR/FPR"
$ws.Cells.Item(2, 25).Formula = "=FNR/TNR"
$ws.Cells.Item(2, 26).Formula = "=TPR+TNR-1"
$ws.Cells.Item(2, 27).Formula = "=(TPR+TNR)/2"
$ws.Cells.Item(2, 28).Formula = "=SQRT(TPR*TNR)"
$ws.Cells.Item(2, 29).Formula = "=PPV+NPV-1"
$ws.Cells.Item(2, 30).Formula = "=(2*PPV*TPR)/(PPV+TPR)"
$ws.Cells.Item(2, 31).Formula = "=(5*PPV*TPR)/(4*PPV+TPR)"
$ws.Cells.Item(2, 32).Formula = "=(1.25*PPV*TPR)/(0.25*PPV+TPR)"
$ws.Cells.Item(2, 33).Formula = "=((P*OP)+(N*ON))/Sn^2"
$ws.Cells.Item(2, 34).Formula = "=(ACC-CKC)/(1-CKC)"
$ws.Cells.Item(2, 35).Formula = "=SQRT(INFORM*MARK)"
$ws.Cells.Item(2, 36).Formula = "=(TP-TN)/(FP-FN)"
$ws.Cells.Item(2, 37).Formula = "=SQRT(3)/PI()*LN(AJ2)"
$ws.Cells.Item(2, 38).Formula = "=P/Sn"
$ws.Cells.Item(2, 39).Formula = "=N/P"

# --- 6. Leftover scratch cells (rows 20-21, cols V/W) on "worksheet" ---
$ws.Cells.Item(20, 22).Value = 1
$ws.Cells.Item(21, 22).Value = 2
$ws.Cells.Item(21, 23).Formula = "=V21+V20"

# --- 7. Populate the new "metrics" sheet (transposed METRIC/VALUE table) ---
$wsMetrics.Cells.Item(1, 1).Value = "METRIC"
$wsMetrics.Cells.Item(1, 2).Value = "VALUE"
$wsMetrics.Cells.Item(2, 1).Value = "Sn"
$wsMetrics.Cells.Item(2, 2).Formula = "=COUNT(y)"
$wsMetrics.Cells.Item(3, 1).Value = "TP"
$wsMetrics.Cells.Item(3, 2).Formula = "=COUNTIFS(y,1,y_pred,1)"
$wsMetrics.Cells.Item(4, 1).Value = "FP"
$wsMetrics.Cells.Item(4, 2).Formula = "=COUNTIFS(y,0,y_pred,1)"
$wsMetrics.Cells.Item(5, 1).Value = "FN"
$wsMetrics.Cells.Item(5, 2).Formula = "=COUNTIFS(y,1,y_pred,0)"
$wsMetrics.Cells.Item(6, 1).Value = "TN"
$wsMetrics.Cells.Item(6, 2).Formula = "=COUNTIFS(y,0,y_pred,0)"
$wsMetrics.Cells.Item(7, 1).Value = "TC"
$wsMetrics.Cells.Item(7, 2).Formula = "=TP+TN"
$wsMetrics.Cells.Item(8, 1).Value = "FC"
$wsMetrics.Cells.Item(8, 2).Formula = "=FP+FN"
$wsMetrics.Cells.Item(9, 1).Value = "OP"
$wsMetrics.Cells.Item(9, 2).Formula = "=TP+FP"
$wsMetrics.Cells.Item(10, 1).Value = "ON"
$wsMetrics.Cells.Item(10, 2).Formula = "=FN+TN"
$wsMetrics.Cells.Item(11, 1).Value = "P"
$wsMetrics.Cells.Item(11, 2).Formula = "=TP+FN"
$wsMetrics.Cells.Item(12, 1).Value = "N"
$wsMetrics.Cells.Item(12, 2).Formula = "=FP+TN"
$wsMetrics.Cells.Item(13, 1).Value = "TPR"
$wsMetrics.Cells.Item(13, 2).Formula = "=TP/P"
$wsMetrics.Cells.Item(14, 1).Value = "FNR"
$wsMetrics.Cells.Item(14, 2).Formula = "=FN/P"
$wsMetrics.Cells.Item(15, 1).Value = "TNR"
$wsMetrics.Cells.Item(15, 2).Formula = "=TN/N"
$wsMetrics.Cells.Item(16, 1).Value = "FPR"
$wsMetrics.Cells.Item(16, 2).Formula = "=FP/N"
$wsMetrics.Cells.Item(17, 1).Value = "ACC"
$wsMetrics.Cells.Item(17, 2).Formula = "=TC/Sn"
$wsMetrics.Cells.Item(18, 1).Value = "DR"
$wsMetrics.Cells.Item(18, 2).Formula = "=TP/Sn"
$wsMetrics.Cells.Item(19, 1).Value = "CRR"
$wsMetrics.Cells.Item(19, 2).Formula = "=TN/Sn"
$wsMetrics.Cells.Item(20, 1).Value = "PPV"
$wsMetrics.Cells.Item(20, 2).Formula = "=TP/OP"
$wsMetrics.Cells.Item(21, 1).Value = "FDR"
$wsMetrics.Cells.Item(21, 2).Formula = "=FP/OP"
$wsMetrics.Cells.Item(22, 1).Value = "FOR"
$wsMetrics.Cells.Item(22, 2).Formula = "=FN/ON"
$wsMetrics.Cells.Item(23, 1).Value = "NPV"
$wsMetrics.Cells.Item(23, 2).Formula = "=TN/ON"
$wsMetrics.Cells.Item(24, 1).Value = "MCR"
$wsMetrics.Cells.Item(24, 2).Formula = "=FC/Sn"
$wsMetrics.Cells.Item(25, 1).Value = "LRP"
$wsMetrics.Cells.Item(25, 2).Formula = "=TPR/FPR"
$wsMetrics.Cells.Item(26, 1).Value = "LRN"
$wsMetrics.Cells.Item(26, 2).Formula = "=FNR/TNR"
$wsMetrics.Cells.Item(27, 1).Value = "INFORM"
$wsMetrics.Cells.Item(27, 2).Formula = "=TPR+TNR-1"
$wsMetrics.Cells.Item(28, 1).Value = "BACC"
$wsMetrics.Cells.Item(28, 2).Formula = "=(TPR+TNR)/2"
$wsMetrics.Cells.Item(29, 1).Value = "GM"
$wsMetrics.Cells.Item(29, 2).Formula = "=SQRT(TPR*TNR)"
$wsMetrics.Cells.Item(30, 1).Value = "MARK"
$wsMetrics.Cells.Item(30, 2).Formula = "=PPV+NPV-1"
$wsMetrics.Cells.Item(31, 1).Value = "F1"
$wsMetrics.Cells.Item(31, 2).Formula = "=(2*PPV*TPR)/(PPV+TPR)"
$wsMetrics.Cells.Item(32, 1).Value = "F2"
$wsMetrics.Cells.Item(32, 2).Formula = "=(5*PPV*TPR)/(4*PPV+TPR)"
$wsMetrics.Cells.Item(33, 1).Value = "F05"
$wsMetrics.Cells.Item(33, 2).Formula = "=(1.25*PPV*TPR)/(0.25*PPV+TPR)"
$wsMetrics.Cells.Item(34, 1).Value = "CKC"
$wsMetrics.Cells.Item(34, 2).Formula = "=((P*OP)+(N*ON))/Sn^2"
$wsMetrics.Cells.Item(35, 1).Value = "CK"
$wsMetrics.Cells.Item(35, 2).Formula = "=(ACC-CKC)/(1-CKC)"
$wsMetrics.Cells.Item(36, 1).Value = "MCC"
$wsMetrics.Cells.Item(36, 2).Formula = "=SQRT(INFORM*MARK)"
$wsMetrics.Cells.Item(37, 1).Value = "OR"
$wsMetrics.Cells.Item(37, 2).Formula = "=(TP-TN)/(FP-FN)"
$wsMetrics.Cells.Item(38, 1).Value = "DP"
$wsMetrics.Cells.Item(38, 2).Value = -0.43469932557681495
$wsMetrics.Cells.Item(39, 1).Value = "PREV"
$wsMetrics.Cells.Item(39, 2).Formula = "=P/Sn"
$wsMetrics.Cells.Item(40, 1).Value = "SKEW"
$wsMetrics.Cells.Item(40, 2).Formula = "=N/P"

# --- 8. View state: active cells / selections / scroll position / active tab ---
$wsData.Activate()
$wsData.Range("A71").Select()
$excel.ActiveWindow.ScrollRow = 63
$excel.ActiveWindow.ScrollColumn = 1

$ws.Activate()
$ws.Range("AK2").Select()
$excel.ActiveWindow.ScrollColumn = 21
$excel.ActiveWindow.ScrollRow = 1

$wsMetrics.Activate()
$wsMetrics.Range("B38").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1

